$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/15/2025  Through  12/21/2025"

# --- Crime Complaints table updates (rows 15-30) ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 34
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 36
$ws.Range("L15").Value = 78.947368421052
$ws.Range("M15").Value = 112.5
$ws.Range("N15").Value = -19.047619047619

# Row 16
$ws.Range("C16").Value = "0"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 97
$ws.Range("J16").Value = 131
$ws.Range("K16").Value = -25.954198473282
$ws.Range("L16").Value = 7.777777777777
$ws.Range("M16").Value = -19.166666666666
$ws.Range("N16").Value = -86.023054755043

# Row 17
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 12.5
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 21.739130434782
$ws.Range("I17").Value = 366
$ws.Range("J17").Value = 345
$ws.Range("K17").Value = 6.086956521739
$ws.Range("L17").Value = 40.76923076923
$ws.Range("M17").Value = 137.662337662338
$ws.Range("N17").Value = -28.793774319066

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 73
$ws.Range("J18").Value = 79
$ws.Range("K18").Value = -7.59493670886
$ws.Range("L18").Value = 2.81690140845
$ws.Range("M18").Value = -45.522388059701
$ws.Range("N18").Value = -90.318302387267

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = -5.555555555555
$ws.Range("I19").Value = 178
$ws.Range("J19").Value = 183
$ws.Range("K19").Value = -2.732240437158
$ws.Range("L19").Value = 1.714285714285
$ws.Range("M19").Value = 81.632653061224
$ws.Range("N19").Value = -44.548286604361

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 40
$ws.Range("J20").Value = 75
$ws.Range("K20").Value = -46.666666666666
$ws.Range("L20").Value = -47.368421052631
$ws.Range("M20").Value = -48.051948051948
$ws.Range("N20").Value = -91.189427312775

# Row 21
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -21.739130434782
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 66
$ws.Range("H21").Value = -7.575757575757
$ws.Range("I21").Value = 791
$ws.Range("J21").Value = 842
$ws.Range("K21").Value = -6.05700712589
$ws.Range("L21").Value = 13.812949640287
$ws.Range("M21").Value = 30.743801652892
$ws.Range("N21").Value = -71.69946332737

# Row 22
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = "0"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = -100
$ws.Range("I22").Value = 12
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = 9.090909090909
$ws.Range("L22").Value = -25
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = "***.*"

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 73
$ws.Range("J23").Value = 83
$ws.Range("K23").Value = -12.048192771084
$ws.Range("L23").Value = -17.045454545454
$ws.Range("M23").Value = 87.179487179487
$ws.Range("N23").Value = "***.*"

# Row 24
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 53.333333333333
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 44.827586206896
$ws.Range("I24").Value = 899
$ws.Range("J24").Value = 664
$ws.Range("K24").Value = 35.39156626506
$ws.Range("L24").Value = 57.167832167832
$ws.Range("M24").Value = 156.125356125356
$ws.Range("N24").Value = "***.*"

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 1300
$ws.Range("I25").Value = 178
$ws.Range("J25").Value = 47
$ws.Range("K25").Value = 278.723404255319
$ws.Range("L25").Value = 109.411764705882
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"

# Row 26
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = 8.695652173913
$ws.Range("I26").Value = 510
$ws.Range("J26").Value = 445
$ws.Range("K26").Value = 14.606741573033
$ws.Range("L26").Value = 24.694376528117
$ws.Range("M26").Value = 27.18204488778
$ws.Range("N26").Value = "***.*"

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = 18.918918918918
$ws.Range("L27").Value = 22.222222222222
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -60
$ws.Range("I28").Value = 39
$ws.Range("J28").Value = 49
$ws.Range("K28").Value = -20.408163265306
$ws.Range("L28").Value = 11.428571428571
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"

# Row 29
$ws.Range("C29").Value = "0"
$ws.Range("D29").Value = "0"
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = "0"
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = -100
$ws.Range("I29").Value = 13
$ws.Range("J29").Value = 18
$ws.Range("K29").Value = -27.777777777777
$ws.Range("L29").Value = -7.142857142857
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -70.454545454545

# Row 30
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = "0"
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
$ws.Range("I30").Value = 11
$ws.Range("J30").Value = 16
$ws.Range("K30").Value = -31.25
$ws.Range("L30").Value = -8.333333333333
$ws.Range("M30").Value = -54.166666666666
$ws.Range("N30").Value = -73.809523809523
